$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell is a text value in the source data (percent/price strings,
# including thousand-dot-formatted numbers). A leading apostrophe forces
# Excel to store the literal text instead of re-parsing look-alike numbers
# as numeric values; resetting the style afterwards drops the quote-prefix
# cell format so the cell keeps its original (unstyled) appearance.

$ws.Range("D2").Value = "'59.443.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.96%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.596.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.48%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'536.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.51%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'140.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.56%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.611.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.27%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.09%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +4.12%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +3.73%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +2.97%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.057.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.20%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'59.347.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.80%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'20.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.20%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.606.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.32%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +2.81%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'346.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.60%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +1.82%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.52%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.91%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.01%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'67.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.58%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.65%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +2.62%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.06%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.36%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'USDe"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.04%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'PEPE"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0738"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.03%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +5.65%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.21%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'18.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.71%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'149.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.07%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.45%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +1.98%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +2.35%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +5.26%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.847"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.834"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.14%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +1.77%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'FirstDigitalUSD"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.02%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Bittensor"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'276.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.65%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.598"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.60%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'10.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.69%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +2.92%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +1.98%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Maker"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1.943.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.28%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'VeChain"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.0223"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.25%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'InjectiveProtocol"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'18.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.46%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'RenderToken"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'4.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.77%  "
$ws.Range("E51").Style = "Normal"
